$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (D) and 1h volume change (E) columns
# with refreshed values from the latest data pull.

$ws.Range("D2").Value = '30.429.04'
$ws.Range("E2").Value = '  -1.09%  '

$ws.Range("D3").Value = '2.091.29'
$ws.Range("E3").Value = '  -1.20%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '330.14'
$ws.Range("E5").Value = '  -1.38%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("E7").Value = '  -0.43%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4374'
$ws.Range("E8").Value = '  -0.80%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.63'
$ws.Range("E9").Value = '  +13.82%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.08876'
$ws.Range("E10").Value = '  -2.33%  '

$ws.Range("E11").Value = '  -2.65%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '24.27'
$ws.Range("E12").Value = '  -4.07%  '

$ws.Range("D13").Value = '2.086.98'
$ws.Range("E13").Value = '  -1.35%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.696'
$ws.Range("E14").Value = '  -1.17%  '

$ws.Range("E15").Value = '  -1.96%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '95.84'
$ws.Range("E16").Value = '  -2.32%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  +0.12%  '

$ws.Range("E18").Value = '  -1.55%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06588'
$ws.Range("E19").Value = '  -0.87%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.19'
$ws.Range("E20").Value = '  -0.16%  '

$ws.Range("E21").Value = '  +0.05%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.259'
$ws.Range("E22").Value = '  -2.29%  '

$ws.Range("D23").Value = '30.467.50'
$ws.Range("E23").Value = '  -1.41%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.24'
$ws.Range("E24").Value = '  +1.26%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.332'
$ws.Range("E25").Value = '  +3.53%  '

$ws.Range("D26").Value = '2.333.55'
$ws.Range("E26").Value = '  -1.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '22.24'
$ws.Range("E27").Value = '  -3.26%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.561'
$ws.Range("E28").Value = '  -0.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '162.70'
$ws.Range("E29").Value = '  -0.46%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '131.48'
$ws.Range("E30").Value = '  -1.70%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.183'
$ws.Range("E31").Value = '  -0.10%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.1068'
$ws.Range("E32").Value = '  -0.55%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.660'
$ws.Range("E33").Value = '  +7.74%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.167'
$ws.Range("E34").Value = '  -1.67%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.895'
$ws.Range("E35").Value = '  -0.66%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '10.03'
$ws.Range("E36").Value = '  +3.94%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02569'
$ws.Range("E37").Value = '  -1.27%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.06825'
$ws.Range("E38").Value = '  +0.65%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.465'
$ws.Range("E39").Value = '  -2.46%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.64'
$ws.Range("E40").Value = '  -1.30%  '

$ws.Range("E41").Value = '  -0.86%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6887'
$ws.Range("E42").Value = '  +0.72%  '

$ws.Range("E43").Value = '  -0.14%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9998'
$ws.Range("E44").Value = '  +0.07%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6344'
$ws.Range("E45").Value = '  -1.61%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.91'
$ws.Range("E46").Value = '  -2.23%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.197'
$ws.Range("E47").Value = '  -3.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.621'
$ws.Range("E48").Value = '  -1.55%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.234'
$ws.Range("E49").Value = '  +7.33%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.243'
$ws.Range("E50").Value = '  -3.45%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '81.79'
$ws.Range("E51").Value = '  -1.74%  '
